# Apply updated odds values to Sheet1 (workbook already open as $excel.ActiveWorkbook)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H3").Value = 2.8
$ws.Range("I3").Value = 3.8
$ws.Range("J3").Value = 3.2
$ws.Range("M3").Value = 1.14
$ws.Range("N3").Value = 5.5
$ws.Range("W3").Value = 2.25
$ws.Range("X3").Value = 1.57
$ws.Range("AF3").Value = 6
$ws.Range("AG3").Value = 21
$ws.Range("AJ3").Value = 7.5

$ws.Range("G5").Value = 2.8
$ws.Range("H5").Value = 2.88
$ws.Range("J5").Value = 3.6
$ws.Range("O5").Value = 1.57
$ws.Range("P5").Value = 2.25
$ws.Range("Q5").Value = 2.88
$ws.Range("R5").Value = 1.4
$ws.Range("S5").Value = 6
$ws.Range("T5").Value = 1.13
$ws.Range("AR5").Value = 4.29
$ws.Range("AS5").Value = 1.2

$ws.Range("G8").Value = 2.2
$ws.Range("I8").Value = 2.8
$ws.Range("L8").Value = 3.4
$ws.Range("U8").Value = 1.3
$ws.Range("V8").Value = 3.4
$ws.Range("Y8").Value = 10
$ws.Range("AB8").Value = 21
$ws.Range("AD8").Value = 21
$ws.Range("AO8").Value = 26

$ws.Range("O9").Value = 1.25
$ws.Range("P9").Value = 3.75
$ws.Range("Q9").Value = 1.88
$ws.Range("R9").Value = 1.98

$ws.Range("N10").Value = 12
$ws.Range("Q10").Value = 1.75
$ws.Range("R10").Value = 2.05
$ws.Range("U10").Value = 1.36
$ws.Range("V10").Value = 3
$ws.Range("AE10").Value = 12
$ws.Range("AN10").Value = 19

$ws.Range("G11").Value = 1.57
$ws.Range("I11").Value = 5.75
$ws.Range("N11").Value = 8.5
$ws.Range("Q11").Value = 2.15
$ws.Range("R11").Value = 1.67
$ws.Range("AB11").Value = 11
$ws.Range("AF11").Value = 7.5
$ws.Range("AJ11").Value = 12
$ws.Range("AL11").Value = 19
$ws.Range("AQ11").Value = 2.24

$ws.Range("H13").Value = 2.72
$ws.Range("I13").Value = 2.95
$ws.Range("J13").Value = 3.4
$ws.Range("L13").Value = 3.65
$ws.Range("P13").Value = 2.32
$ws.Range("Q13").Value = 2.57
$ws.Range("V13").Value = 2.25
$ws.Range("W13").Value = 2.05
$ws.Range("Y13").Value = 6.3
$ws.Range("Z13").Value = 11.75
$ws.Range("AA13").Value = 10.5
$ws.Range("AC13").Value = 28
$ws.Range("AH13").Value = 110
$ws.Range("AJ13").Value = 6.9
$ws.Range("AL13").Value = 11

$ws.Range("Q14").Value = 2.5
$ws.Range("R14").Value = 1.53
$ws.Range("AP14").Value = 1.83
$ws.Range("AQ14").Value = 2.03

$ws.Range("M15").Value = 1.03
$ws.Range("N15").Value = 15
$ws.Range("Q15").Value = 1.73
$ws.Range("R15").Value = 2.08

$ws.Range("G16").Value = 2.4
$ws.Range("I16").Value = 2.75
$ws.Range("L16").Value = 3.4
$ws.Range("M16").Value = 1.06
$ws.Range("N16").Value = 10
$ws.Range("O16").Value = 1.3
$ws.Range("P16").Value = 3.4
$ws.Range("Q16").Value = 2
$ws.Range("R16").Value = 1.85
$ws.Range("AC16").Value = 21
$ws.Range("AE16").Value = 10
$ws.Range("AK16").Value = 13

$ws.Range("G17").Value = 1.53
$ws.Range("I17").Value = 6
$ws.Range("K17").Value = 2.3
$ws.Range("U17").Value = 1.36
$ws.Range("V17").Value = 3
$ws.Range("Y17").Value = 7
$ws.Range("AD17").Value = 26
$ws.Range("AE17").Value = 11
$ws.Range("AG17").Value = 17
$ws.Range("AI17").Value = 301
$ws.Range("AK17").Value = 34
$ws.Range("AN17").Value = 41

$ws.Range("M19").Value = 1.06
$ws.Range("N19").Value = 10
$ws.Range("Q19").Value = 2.05
$ws.Range("R19").Value = 1.75

$ws.Range("G21").Value = 2.38
$ws.Range("I21").Value = 2.88
$ws.Range("AJ21").Value = 9.5

$ws.Range("J24").Value = 2.63
$ws.Range("AB24").Value = 15
$ws.Range("AE24").Value = 7.5
$ws.Range("AM24").Value = 51

$ws.Range("G26").Value = 3

$ws.Range("M27").Value = 1.05
$ws.Range("N27").Value = 11
$ws.Range("O27").Value = 1.29
$ws.Range("P27").Value = 3.5
$ws.Range("Q27").Value = 1.9
$ws.Range("R27").Value = 1.9
$ws.Range("S27").Value = 3.25
$ws.Range("T27").Value = 1.33

$ws.Range("I28").Value = 3
$ws.Range("AE28").Value = 11.25
$ws.Range("AG28").Value = 11.75
$ws.Range("AI28").Value = 250

$ws.Range("M29").Value = 1.04
$ws.Range("O29").Value = 1.22
$ws.Range("T29").Value = 1.4

$ws.Range("H30").Value = 3.3
$ws.Range("I30").Value = 2.88
$ws.Range("J30").Value = 3.2
$ws.Range("M30").Value = 1.07
$ws.Range("N30").Value = 9
$ws.Range("O30").Value = 1.33
$ws.Range("P30").Value = 3.25
$ws.Range("T30").Value = 1.29
$ws.Range("AE30").Value = 9
$ws.Range("AG30").Value = 15
$ws.Range("AK30").Value = 13

$ws.Range("M31").Value = 1.1
$ws.Range("O31").Value = 1.44
$ws.Range("P31").Value = 2.63
$ws.Range("T31").Value = 1.18

$ws.Range("AD32").Value = 23
$ws.Range("AO32").Value = 23

$ws.Range("G33").Value = 2.7
$ws.Range("I33").Value = 2.5
$ws.Range("J33").Value = 3.25
$ws.Range("L33").Value = 3.1
$ws.Range("M33").Value = 1.06
$ws.Range("N33").Value = 8
$ws.Range("Q33").Value = 2.05
$ws.Range("R33").Value = 1.8
$ws.Range("S33").Value = 3.5
$ws.Range("T33").Value = 1.29
$ws.Range("AD33").Value = 34
$ws.Range("AJ33").Value = 8.5
$ws.Range("AL33").Value = 10
$ws.Range("AM33").Value = 23
